$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateStyleSrc = $ws.Cells.Item(2, 1)

$dates = @(43492, 43493, 43494, 43495, 43496, 43497, 43498)

$row = 3
foreach ($d in $dates) {
    $dateStyleSrc.Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = $d
    $row++
}

# Text values are entered in this specific order so that the shared-string
# table is built up exactly like the original authoring session.
$ws.Cells.Item(3, 2).Value = "no updates"
$ws.Cells.Item(5, 2).Value = "Eclipse setup"
$ws.Cells.Item(4, 2).Value = "meeting on task updates"
$ws.Cells.Item(6, 2).Value = "java and spring basics"
$ws.Cells.Item(7, 2).Value = "started the task"
$ws.Cells.Item(8, 2).Value = "no updates"
$ws.Cells.Item(9, 2).Value = "continued the task"

$ws.Range("A10").Select()
